$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Npc")

# Rows 12-21: Cid block 13100-13109 -> 14100-14109
$ws.Range("A12").Value = 14100
$ws.Range("A13").Value = 14101
$ws.Range("A14").Value = 14102
$ws.Range("A15").Value = 14103
$ws.Range("A16").Value = 14104
$ws.Range("A17").Value = 14105
$ws.Range("A18").Value = 14106
$ws.Range("A19").Value = 14107
$ws.Range("A20").Value = 14108
$ws.Range("A21").Value = 14109

# Rows 52-63: Cid block 14100-14111 -> 13100-13111
$ws.Range("A52").Value = 13100
$ws.Range("A53").Value = 13101
$ws.Range("A54").Value = 13102
$ws.Range("A55").Value = 13103
$ws.Range("A56").Value = 13104
$ws.Range("A57").Value = 13105
$ws.Range("A58").Value = 13106
$ws.Range("A59").Value = 13107
$ws.Range("A60").Value = 13108
$ws.Range("A61").Value = 13109
$ws.Range("A62").Value = 13110
$ws.Range("A63").Value = 13111

# Rows 64-69: Cid block 13100-13105 -> 12100-12105
$ws.Range("A64").Value = 12100
$ws.Range("A65").Value = 12101
$ws.Range("A66").Value = 12102
$ws.Range("A67").Value = 12103
$ws.Range("A68").Value = 12104
$ws.Range("A69").Value = 12105

# Rows 70-74: Cid block 12100-12103,12103 -> 11100-11104
$ws.Range("A70").Value = 11100
$ws.Range("A71").Value = 11101
$ws.Range("A72").Value = 11102
$ws.Range("A73").Value = 11103
$ws.Range("A74").Value = 11104

# Update the view: select A74 and scroll so row 31 is the top visible row
$ws.Activate()
$ws.Range("A74").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
